# Restore C10 ("From" value for rule R30) from 18 to 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
